$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E34").Value = 3
$ws.Range("F34").Value = "(-4)For 4 test cases failed."
